$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text format so numeric-looking
# strings (e.g. '13.00', '6.70') are not coerced into numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '47.275.25'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +4.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.496.12'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.89%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.08'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.93'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.76%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.527'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.44%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.32'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +7.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0815'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.09%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.53'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.22'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.882.92'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.504.56'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.858'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '47.205.62'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.00'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +6.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.70'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +5.59%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.83'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.47'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +8.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '250.88'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.49%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.28'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.13%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.30'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.05'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.04'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +6.37%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +10.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.46'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Celestia'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.73'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.41%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.49'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +5.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0794'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.24%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.99'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +6.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.72'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.48%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.91%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '122.43'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.40%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.19'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0299'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.969.44'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.02'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.48%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.53%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.06'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.70%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +8.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.44'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.51%  '
